# Auto-generated edit script: refresh cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain TEXT even though it looks numeric
# (e.g. "150.80"), without permanently altering the target cell's number
# format/style. We stage the text in a scratch cell (far outside the used
# range), force it to Text format there, copy/paste-special VALUES ONLY into
# the destination (so only the string payload moves, not the scratch cell's
# style), then clear the scratch cell so no trace / used-range growth remains.
function Set-TextValue {
    param([string]$CellRef, [string]$Val)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $Val
    $scratch.Copy()
    $dest = $ws.Range($CellRef)
    $dest.PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

$excel.CutCopyMode = $false

$ws.Range("D2").Value = "57.282.75"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.426.98"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue "D5" "488.94"
$ws.Range("E5").Value = "  -0.15%  "
Set-TextValue "D6" "154.74"
$ws.Range("E6").Value = "  +2.06%  "
Set-TextValue "D7" "0.617"
$ws.Range("E7").Value = "  +19.47%  "
Set-TextValue "D8" "0.997"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "2.445.56"
$ws.Range("E9").Value = "  -1.29%  "
Set-TextValue "D10" "6.18"
$ws.Range("E10").Value = "  +8.18%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.845.96"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "57.243.29"
$ws.Range("E15").Value = "  +0.14%  "
Set-TextValue "D16" "20.61"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "2.443.48"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +1.50%  "
Set-TextValue "D20" "324.53"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  -1.87%  "
Set-TextValue "D22" "0.997"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.93%  "
Set-TextValue "D24" "57.84"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D25" "0.401"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D26" "0.994"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "2.534.63"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").Value = "0.0₃0787"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("E31").Value = "  -0.02%  "
Set-TextValue "D32" "150.80"
$ws.Range("E32").Value = "  -0.02%  "
Set-TextValue "D33" "18.70"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("E34").Value = "  +0.31%  "
Set-TextValue "D35" "5.30"
$ws.Range("E35").Value = "  +1.75%  "
Set-TextValue "D36" "3.79"
$ws.Range("E36").Value = "  +0.64%  "
Set-TextValue "D37" "1.14"
$ws.Range("E37").Value = "  -0.71%  "
Set-TextValue "D38" "0.819"
$ws.Range("E38").Value = "  -7.93%  "
$ws.Range("E39").Value = "  +7.59%  "
Set-TextValue "D40" "285.48"
$ws.Range("E40").Value = "  +8.38%  "
Set-TextValue "D41" "34.05"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "1.38"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "3.52"
$ws.Range("E43").Value = "  +0.48%  "
Set-TextValue "D44" "0.993"
$ws.Range("E44").Value = "  -0.29%  "
Set-TextValue "D45" "0.602"
$ws.Range("E45").Value = "  -0.92%  "
Set-TextValue "D46" "0.0532"
$ws.Range("E46").Value = "  -4.59%  "
$ws.Range("E47").Value = "  -0.06%  "
Set-TextValue "D48" "0.0228"
$ws.Range("E48").Value = "  -0.38%  "
Set-TextValue "D49" "4.57"
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("D50").Value = "1.898.22"
$ws.Range("E50").Value = "  +2.15%  "
Set-TextValue "D51" "17.61"

$excel.CutCopyMode = $false
Write-Host "Done applying cryptos refresh."
